# Actualizacion desde MV -datos-
# Append five new daily rows to the bottom of the data table, following the
# same "no auction held" pattern already used by several existing rows
# (Serie/date in col A, Cupo = 10000 in col B, 0 in col D, C/E/F/G left
# blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @("08-09-2021", "09-09-2021", "14-09-2021", "15-09-2021", "16-09-2021")

$startRow = 20
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = $startRow + $i

    # The date strings look like valid dates, so a plain .Value assignment
    # would be auto-converted to a date serial by Excel's type inference.
    # Prefixing with an apostrophe forces it to be entered as text (same
    # as the pre-existing date-label cells in column A), and ClearFormats
    # drops the quote-prefix cell style back to the sheet's default so the
    # new cells end up styled exactly like the other plain data rows.
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = "'" + $newDates[$i]
    $cell.ClearFormats()

    $ws.Cells.Item($r, 2).Value = 10000
    $ws.Cells.Item($r, 4).Value = 0
}
